# Commit: "Evaluate HS for different cases / And unmet demand / wasted surplus"
#
# This script reproduces the data + view edits made to
# model/Inputs/model_inputs_inelas.xlsx:
#   - parameters: "min SoC" (A2) 0.2 -> 0.02
#   - tech: Owned PV UOVC (G3) 1.38E-2 -> 0
#           Owned Batteries Lifetime (D4) 12 -> 8
#           Owned Batteries UOVC (G4) 0.6 -> 0.0006, shown with 4 decimals
#   - day_weights: Day 1/2/3 weights (B2:B4) 91/153/121 -> 199/106/61
#   - cap_factors: full recompute of the summer/fall_spring/winter capacity
#     factor curves (rows 2-4, columns G:U)
#   - final selections / active sheet match the saved workbook state
#     (day_weights ends up the active tab)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# parameters sheet: min SoC 0.2 -> 0.02
# ---------------------------------------------------------------------
$wsParams = $wb.Worksheets.Item("parameters")
$wsParams.Range("A2").Value = 0.02

# ---------------------------------------------------------------------
# tech sheet: Owned PV UOVC, Owned Batteries Lifetime + UOVC
# ---------------------------------------------------------------------
$wsTech = $wb.Worksheets.Item("tech")
$wsTech.Range("G3").Value = 0
$wsTech.Range("D4").Value = 8
$wsTech.Range("G4").Value = 0.0006
$wsTech.Range("G4").NumberFormat = "0.0000"
# Column G widened slightly after the format change (best-effort match of
# Excel's own auto best-fit recompute for the new cell content)
$wsTech.Columns.Item(7).ColumnWidth = 5.5

# ---------------------------------------------------------------------
# day_weights sheet: Day 1/2/3 weights
# ---------------------------------------------------------------------
$wsDays = $wb.Worksheets.Item("day_weights")
$wsDays.Range("B2").Value = 199
$wsDays.Range("B3").Value = 106
$wsDays.Range("B4").Value = 61

# ---------------------------------------------------------------------
# cap_factors sheet: recomputed capacity-factor curves (rows 2-4, G:U)
# ---------------------------------------------------------------------
$wsCap = $wb.Worksheets.Item("cap_factors")
$wsCap.Range("G2").Value = 0.00402010050251256
$wsCap.Range("H2").Value = 0.04117085427135676
$wsCap.Range("I2").Value = 0.1549849246231156
$wsCap.Range("J2").Value = 0.3412010050251257
$wsCap.Range("K2").Value = 0.5113718592964824
$wsCap.Range("L2").Value = 0.6344020100502513
$wsCap.Range("M2").Value = 0.705894472361809
$wsCap.Range("N2").Value = 0.7173165829145729
$wsCap.Range("O2").Value = 0.6773417085427136
$wsCap.Range("P2").Value = 0.5780954773869347
$wsCap.Range("Q2").Value = 0.4341256281407035
$wsCap.Range("R2").Value = 0.2525276381909548
$wsCap.Range("S2").Value = 0.09160804020100505
$wsCap.Range("T2").Value = 0.01759798994974874
$wsCap.Range("U2").Value = 0.0003216080402010047
$wsCap.Range("G3").Value = 0.002150943396226415
$wsCap.Range("H3").Value = 0.02077358490566038
$wsCap.Range("I3").Value = 0.07405660377358492
$wsCap.Range("J3").Value = 0.2032169811320755
$wsCap.Range("K3").Value = 0.3473301886792453
$wsCap.Range("L3").Value = 0.4670283018867925
$wsCap.Range("M3").Value = 0.5449999999999999
$wsCap.Range("N3").Value = 0.5594245283018868
$wsCap.Range("O3").Value = 0.5117641509433962
$wsCap.Range("P3").Value = 0.4107358490566038
$wsCap.Range("Q3").Value = 0.2745660377358491
$wsCap.Range("R3").Value = 0.1335471698113208
$wsCap.Range("S3").Value = 0.03995283018867922
$wsCap.Range("T3").Value = 0.008669811320754714
$wsCap.Range("U3").Value = 0.0001415094339622642
$wsCap.Range("G4").Value = 0.00180327868852459
$wsCap.Range("H4").Value = 0.01603278688524591
$wsCap.Range("I4").Value = 0.05422950819672132
$wsCap.Range("J4").Value = 0.1274918032786885
$wsCap.Range("K4").Value = 0.2095081967213114
$wsCap.Range("L4").Value = 0.2680983606557377
$wsCap.Range("M4").Value = 0.2826721311475409
$wsCap.Range("N4").Value = 0.2759508196721313
$wsCap.Range("O4").Value = 0.2548852459016394
$wsCap.Range("P4").Value = 0.2180983606557376
$wsCap.Range("Q4").Value = 0.151
$wsCap.Range("R4").Value = 0.08108196721311473
$wsCap.Range("S4").Value = 0.03142622950819671
$wsCap.Range("T4").Value = 0.007704918032786882
$wsCap.Range("U4").Value = 0.00009836065573770502

# ---------------------------------------------------------------------
# Final selections on each touched sheet (matches saved cursor positions)
# ---------------------------------------------------------------------
$wsParams.Range("A3").Select()
$wsTech.Range("L9").Select()
$wsCap.Range("B3:Y3").Select()

# day_weights is selected/activated last so it ends up the active tab
$wsDays.Range("B2:B4").Select()
$wsDays.Activate()

# Best-effort: scroll the sheet tab strip so "tech" (index 2) is the
# first visible tab, matching the saved window state.
try {
    $excel.ActiveWindow.ScrollWorkbookTabs(3, 3)
} catch {
}
